# Regenerate merged AHB files
# - Rename the "_old"/"_new" suffixed header labels to "_FV2210"/"_FV2304"
# - Wrap the data range in an Excel Table ("Table1")
# - Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row labels (row 1) in place.
$headerRenames = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# 2) Turn the used range into a proper Excel Table ("Table1") with headers.
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3) Freeze the top (header) row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
